# Generate Report for Handback
# The "8ac2b10a-90e9-4c00-a003-c701f467e106.md" file has finished its
# handback cycle: it moved from "Ready for handoff" to
# "Handed back: in sync with en-US" on both locales, the handback
# timestamps advance, and the stale "handback file is not latest" error
# is cleared now that the file is in sync.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-08-29 02:48:15"
$zhcn.Range("P3").Value = ""
# Error Detail column no longer holds the long error message -> narrower
$zhcn.Columns.Item(16).AutoFit()

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-08-29 02:48:23"
$dede.Range("P3").Value = ""
$dede.Columns.Item(16).AutoFit()
